$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-6 (Timestamp in col A, Raw JSON in col J); columns B-I remain "N/A"
$ws.Cells.Item(2, 1).Value = '2025-11-25T04:53:21.081Z'
$ws.Cells.Item(2, 10).Value = '{"Sequence":1,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"9cbf6af0-76d3-43e5-bd6d-e4c6b08283b7","EventDtm":"2025-11-25T04:50:39Z","AppDtm":"2025-11-25T04:51:09Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T04:50:39Z","DeviceID":"JSGA623040193","GPSLockState":"LOCKED","GPSSatelliteCount":10,"GPSLastLock":0,"GPSLatitude":28.678724,"GPSLongitude":77.599637,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"97","LAC":"1827","CID":"230077975","RSSI":"-65","ExtPower":true,"ExtPowerVoltage":30.9,"BatteryVoltage":8,"DeviceTemp":25,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"SJKU4000017","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T04:50:39Z","AssetType":"Reefer","AssetID":"SJKU4000017","OEM":"CARRIER","TAmb":19.82,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":0,"PctCO2Q":null,"PctCO2Set":5,"PctCO2SetQ":null,"PSuc":62.3,"PSucQ":null,"TDis":43.34,"TDisQ":null,"FreqComp":null,"TSuc":17.57,"TSucQ":null,"MCond":0,"PCond":276.17,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Idle","SnCtrl":"04475864","AmpPhA":1.58,"AmpPhB":1.53,"AmpPhC":1.58,"PDis":84.55,"PDisQ":null,"PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":17.64,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":69.55,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":411.81,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0,"PctO2Q":"asProvided","PctO2Set":10,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":17.21,"TRtn1Q":null,"TRtn2":17.21,"TRtn2Q":null,"TSet":17.2,"TSetQ":null,"VerSWMajor":"5370","VerSWMinor":null,"PctSucVlv":20,"PctSucVlvQ":null,"TSup1":17.58,"TSup1Q":null,"TSup2":17.6,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-03-23T08:47:13Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

$ws.Cells.Item(3, 1).Value = '2025-11-25T04:53:48.274Z'
$ws.Cells.Item(3, 10).Value = '{"Sequence":2,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"584a3c61-46c6-4a37-bdd5-b29c92b54789","EventDtm":"2025-11-25T04:53:30Z","AppDtm":"2025-11-25T04:53:47Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T04:53:30Z","DeviceID":"JSGA622340343","GPSLockState":"LOCKED","GPSSatelliteCount":12,"GPSLastLock":0,"GPSLatitude":12.804897,"GPSLongitude":77.662165,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9003","CID":"46904076","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.9,"BatteryVoltage":8,"DeviceTemp":31,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CCLU1035976","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T04:53:30Z","AssetType":"Reefer","AssetID":"CCLU1035976","OEM":"CARRIER","TAmb":31.78,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":5,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":98.72,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04904656","AmpPhA":11.46,"AmpPhB":10.92,"AmpPhC":11.46,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":6.28,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":415.79,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0.1,"PctO2Q":"unknown","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.86,"TRtn1Q":null,"TRtn2":5.88,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"5180","VerSWMinor":null,"PctSucVlv":4.76,"PctSucVlvQ":null,"TSup1":5.1,"TSup1Q":null,"TSup2":5.05,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-05-12T10:01:32Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

$ws.Cells.Item(4, 1).Value = '2025-11-25T04:55:20.042Z'
$ws.Cells.Item(4, 10).Value = '{"Sequence":3,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"936fa6b6-851d-4513-ba3d-056a77c0e455","EventDtm":"2025-11-25T04:55:02Z","AppDtm":"2025-11-25T04:55:19Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T04:55:02Z","DeviceID":"JSGA623040182","GPSLockState":"LOCKED","GPSSatelliteCount":15,"GPSLastLock":0,"GPSLatitude":17.538987,"GPSLongitude":78.17834,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19397","CID":"232155661","RSSI":"-55","ExtPower":true,"ExtPowerVoltage":31.4,"BatteryVoltage":8,"DeviceTemp":46,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TDRU5800414","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-25T04:55:02Z","AssetType":"Reefer","AssetID":"TDRU5800414","OEM":"DAIKIN","TAmb":29.81,"TAmbQ":null,"TUSDA4":-46.4,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":70,"PSucQ":"asProvided","TDis":80.8,"TDisQ":null,"FreqComp":null,"TSuc":-13.9,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":990,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":18,"PctExpVlvQ":null,"TEvap":3.31,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":98.82,"PctHumQ":null,"PctHumSet":75,"PctHumSetQ":"configured","FreqLine":50,"FreqLineQ":null,"VLine1":405.7,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":5.12,"TRtn1Q":null,"TRtn2":5.1,"TRtn2Q":null,"TSet":4,"TSetQ":null,"VerSWMajor":"24C2","VerSWMinor":null,"PctSucVlv":94.21,"PctSucVlvQ":null,"TSup1":4.12,"TSup1Q":null,"TSup2":4.2,"TSup2Q":null,"AmpTotalDraw":13,"AmpTotalDrawQ":"asProvided","TUSDA1":-46.4,"TUSDA1Q":"OOR","TUSDA2":-46.4,"TUSDA2Q":"OOR","TUSDA3":-46.3,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"OOR","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-09-30T07:54:44Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'

$ws.Cells.Item(5, 1).Value = '2025-11-25T04:56:43.542Z'
$ws.Cells.Item(5, 10).Value = '{"Sequence":4,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"c29c27ca-bc62-4aa6-a7e4-9fc140faaa7d","EventDtm":"2025-11-25T04:56:25Z","AppDtm":"2025-11-25T04:56:42Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T04:56:25Z","DeviceID":"JSGA623040311","GPSLockState":"LOCKED","GPSSatelliteCount":11,"GPSLastLock":0,"GPSLatitude":26.606155,"GPSLongitude":80.723748,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"15","LAC":"55363","CID":"210193709","RSSI":null,"ExtPower":true,"ExtPowerVoltage":34.1,"BatteryVoltage":8,"DeviceTemp":23,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":null,"LastAssetRunState":null},"ReeferData":null}}'

$ws.Cells.Item(6, 1).Value = '2025-11-25T04:57:18.149Z'
$ws.Cells.Item(6, 10).Value = '{"Sequence":5,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"b805d92d-ac79-4cc0-a4af-de818c839db0","EventDtm":"2025-11-25T04:57:02Z","AppDtm":"2025-11-25T04:57:17Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-25T04:57:02Z","DeviceID":"JSGA623040277","GPSLockState":"LOCKED","GPSSatelliteCount":13,"GPSLastLock":0,"GPSLatitude":17.572218,"GPSLongitude":78.515077,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"19327","CID":"235323415","RSSI":"-71","ExtPower":true,"ExtPowerVoltage":32.5,"BatteryVoltage":8,"DeviceTemp":29,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"DFOU6120281","LastAssetRunState":"Running"},"ReeferData":null}}'

# Remove the old trailing rows 7-9 (new data set only has 5 rows of alerts)
$ws.Range("A7:J9").EntireRow.Delete()
